$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.31778466666666
$ws.Range("H2").Value = 135.953354
$ws.Range("I2").Value = 0.1102361023838286
$ws.Range("J2").Value = 0.1102361023838286
$ws.Range("M2").Value = 16.05260533333333
$ws.Range("N2").Value = 48.157816
$ws.Range("O2").Value = 0.1752915379534001
$ws.Range("P2").Value = 0.1752915379534001
$ws.Range("Q2").Value = 727.4685118349847
$ws.Range("R2").Value = 6547.216606514863
$ws.Range("S2").Value = 0.01932345592484979
$ws.Range("T2").Value = 0.01932345592484979
$ws.Range("G3").Value = 45.31778466666666
$ws.Range("H3").Value = 135.953354
$ws.Range("I3").Value = 0.1102361023838286
$ws.Range("J3").Value = 0.1102361023838286
$ws.Range("O3").Value = 0.07888758308485012
$ws.Range("P3").Value = 0.07888758308485012
$ws.Range("Q3").Value = 327.3873533145157
$ws.Range("R3").Value = 2946.486179830642
$ws.Range("S3").Value = 0.008696259685754326
$ws.Range("T3").Value = 0.008696259685754326
$ws.Range("G4").Value = 45.31778466666666
$ws.Range("H4").Value = 135.953354
$ws.Range("I4").Value = 0.1102361023838286
$ws.Range("J4").Value = 0.1102361023838286
$ws.Range("M4").Value = 2.098187333333334
$ws.Range("N4").Value = 6.294562000000001
$ws.Range("O4").Value = 0.02291182502385553
$ws.Range("P4").Value = 0.02291182502385553
$ws.Range("Q4").Value = 95.08520176232756
$ws.Range("R4").Value = 855.7668158609481
$ws.Range("S4").Value = 0.002525710289130105
$ws.Range("T4").Value = 0.002525710289130105
$ws.Range("G5").Value = 45.31778466666666
$ws.Range("H5").Value = 135.953354
$ws.Range("I5").Value = 0.1102361023838286
$ws.Range("J5").Value = 0.1102361023838286
$ws.Range("M5").Value = 66.20156266666667
$ws.Range("N5").Value = 198.604688
$ws.Range("O5").Value = 0.7229090539378943
$ws.Range("P5").Value = 0.7229090539378942
$ws.Range("Q5").Value = 3000.108161524839
$ws.Range("R5").Value = 27000.97345372355
$ws.Range("S5").Value = 0.07969067648409442
$ws.Range("T5").Value = 0.07969067648409441
$ws.Range("I6").Value = 0.2429203181515272
$ws.Range("J6").Value = 0.2429203181515272
$ws.Range("M6").Value = 16.05260533333333
$ws.Range("N6").Value = 48.157816
$ws.Range("O6").Value = 0.1752915379534001
$ws.Range("P6").Value = 0.1752915379534001
$ws.Range("Q6").Value = 1603.076292781705
$ws.Range("R6").Value = 14427.68663503535
$ws.Range("S6").Value = 0.04258187616891045
$ws.Range("T6").Value = 0.04258187616891046
$ws.Range("I7").Value = 0.2429203181515272
$ws.Range("J7").Value = 0.2429203181515272
$ws.Range("O7").Value = 0.07888758308485012
$ws.Range("P7").Value = 0.07888758308485012
$ws.Range("S7").Value = 0.01916339678117683
$ws.Range("T7").Value = 0.01916339678117683
$ws.Range("I8").Value = 0.2429203181515272
$ws.Range("J8").Value = 0.2429203181515272
$ws.Range("M8").Value = 2.098187333333334
$ws.Range("N8").Value = 6.294562000000001
$ws.Range("O8").Value = 0.02291182502385553
$ws.Range("P8").Value = 0.02291182502385553
$ws.Range("Q8").Value = 209.5332378786571
$ws.Range("R8").Value = 1885.799140907914
$ws.Range("S8").Value = 0.005565747824227108
$ws.Range("T8").Value = 0.005565747824227108
$ws.Range("I9").Value = 0.2429203181515272
$ws.Range("J9").Value = 0.2429203181515272
$ws.Range("M9").Value = 66.20156266666667
$ws.Range("N9").Value = 198.604688
$ws.Range("O9").Value = 0.7229090539378943
$ws.Range("P9").Value = 0.7229090539378942
$ws.Range("Q9").Value = 6611.148374504926
$ws.Range("R9").Value = 59500.33537054434
$ws.Range("S9").Value = 0.1756092973772128
$ws.Range("T9").Value = 0.1756092973772128
$ws.Range("G10").Value = 16.49037766666667
$ws.Range("H10").Value = 49.471133
$ws.Range("I10").Value = 0.04011305879538658
$ws.Range("J10").Value = 0.04011305879538658
$ws.Range("M10").Value = 16.05260533333333
$ws.Range("N10").Value = 48.157816
$ws.Range("O10").Value = 0.1752915379534001
$ws.Range("P10").Value = 0.1752915379534001
$ws.Range("Q10").Value = 264.7135244806142
$ws.Range("R10").Value = 2382.421720325528
$ws.Range("S10").Value = 0.007031479768258474
$ws.Range("T10").Value = 0.007031479768258474
$ws.Range("G11").Value = 16.49037766666667
$ws.Range("H11").Value = 49.471133
$ws.Range("I11").Value = 0.04011305879538658
$ws.Range("J11").Value = 0.04011305879538658
$ws.Range("O11").Value = 0.07888758308485012
$ws.Range("P11").Value = 0.07888758308485012
$ws.Range("Q11").Value = 119.1307372846454
$ws.Range("R11").Value = 1072.176635561809
$ws.Range("S11").Value = 0.003164422258508537
$ws.Range("T11").Value = 0.003164422258508537
$ws.Range("G12").Value = 16.49037766666667
$ws.Range("H12").Value = 49.471133
$ws.Range("I12").Value = 0.04011305879538658
$ws.Range("J12").Value = 0.04011305879538658
$ws.Range("M12").Value = 2.098187333333334
$ws.Range("N12").Value = 6.294562000000001
$ws.Range("O12").Value = 0.02291182502385553
$ws.Range("P12").Value = 0.02291182502385553
$ws.Range("Q12").Value = 34.5999015420829
$ws.Range("R12").Value = 311.399113878746
$ws.Range("S12").Value = 0.0009190633842915262
$ws.Range("T12").Value = 0.0009190633842915262
$ws.Range("G13").Value = 16.49037766666667
$ws.Range("H13").Value = 49.471133
$ws.Range("I13").Value = 0.04011305879538658
$ws.Range("J13").Value = 0.04011305879538658
$ws.Range("M13").Value = 66.20156266666667
$ws.Range("N13").Value = 198.604688
$ws.Range("O13").Value = 0.7229090539378943
$ws.Range("P13").Value = 0.7229090539378942
$ws.Range("Q13").Value = 1091.688770496834
$ws.Range("R13").Value = 9825.198934471506
$ws.Range("S13").Value = 0.02899809338432804
$ws.Range("T13").Value = 0.02899809338432803
$ws.Range("G14").Value = 249.4253923333333
$ws.Range("H14").Value = 748.276177
$ws.Range("I14").Value = 0.6067305206692575
$ws.Range("J14").Value = 0.6067305206692575
$ws.Range("M14").Value = 16.05260533333333
$ws.Range("N14").Value = 48.157816
$ws.Range("O14").Value = 0.1752915379534001
$ws.Range("P14").Value = 0.1752915379534001
$ws.Range("Q14").Value = 4003.927383238825
$ws.Range("R14").Value = 36035.34644914943
$ws.Range("S14").Value = 0.1063547260913813
$ws.Range("T14").Value = 0.1063547260913813
$ws.Range("G15").Value = 249.4253923333333
$ws.Range("H15").Value = 748.276177
$ws.Range("I15").Value = 0.6067305206692575
$ws.Range("J15").Value = 0.6067305206692575
$ws.Range("O15").Value = 0.07888758308485012
$ws.Range("P15").Value = 0.07888758308485012
$ws.Range("Q15").Value = 1801.913302825424
$ws.Range("R15").Value = 16217.21972542882
$ws.Range("S15").Value = 0.04786350435941043
$ws.Range("T15").Value = 0.04786350435941043
$ws.Range("G16").Value = 249.4253923333333
$ws.Range("H16").Value = 748.276177
$ws.Range("I16").Value = 0.6067305206692575
$ws.Range("J16").Value = 0.6067305206692575
$ws.Range("M16").Value = 2.098187333333334
$ws.Range("N16").Value = 6.294562000000001
$ws.Range("O16").Value = 0.02291182502385553
$ws.Range("P16").Value = 0.02291182502385553
$ws.Range("Q16").Value = 523.3411988054971
$ws.Range("R16").Value = 4710.070789249475
$ws.Range("S16").Value = 0.01390130352620679
$ws.Range("T16").Value = 0.01390130352620679
$ws.Range("G17").Value = 249.4253923333333
$ws.Range("H17").Value = 748.276177
$ws.Range("I17").Value = 0.6067305206692575
$ws.Range("J17").Value = 0.6067305206692575
$ws.Range("M17").Value = 66.20156266666667
$ws.Range("N17").Value = 198.604688
$ws.Range("O17").Value = 0.7229090539378943
$ws.Range("P17").Value = 0.7229090539378943
$ws.Range("Q17").Value = 16512.35074121309
$ws.Range("R17").Value = 148611.1566709178
$ws.Range("S17").Value = 0.438610986692259
$ws.Range("T17").Value = 0.4386109866922589
